$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 542, shifting existing rows 542:583 down to 543:584
$ws.Rows.Item(542).Insert()

# Populate the newly inserted row 542 with the new record's data.
$ws.Cells.Item(542, 1).Value = 7
$ws.Cells.Item(542, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(542, 3).Value = "Ñuble"
$ws.Cells.Item(542, 4).Value = 45013
$ws.Cells.Item(542, 5).Value = 16
$ws.Cells.Item(542, 6).Value = 100114014
$ws.Cells.Item(542, 7).Value = "Betarraga"
$ws.Cells.Item(542, 8).Value = "Sin especificar"
$ws.Cells.Item(542, 9).Value = "Primera"
$ws.Cells.Item(542, 10).Value = 90
$ws.Cells.Item(542, 11).Value = 800
$ws.Cells.Item(542, 12).Value = 1000
$ws.Cells.Item(542, 13).Value = 911
$ws.Cells.Item(542, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(542, 15).Value = "Región del Maule"
$ws.Cells.Item(542, 16).Value = 182
$ws.Cells.Item(542, 17).Value = 5
$ws.Cells.Item(542, 18).Value = "Hortaliza"
